$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '66.613.65'
Set-TextValue $ws.Range("E2") '  -0.18%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.239.15'
Set-TextValue $ws.Range("E3") '  +1.30%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  -0.01%  '

# Row 5
Set-TextValue $ws.Range("D5") '603.37'
Set-TextValue $ws.Range("E5") '  +0.40%  '

# Row 6
Set-TextValue $ws.Range("D6") '156.45'
Set-TextValue $ws.Range("E6") '  -0.55%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  -0.01%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.238.88'
Set-TextValue $ws.Range("E8") '  +1.23%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.546'
Set-TextValue $ws.Range("E9") '  -1.13%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.162'
Set-TextValue $ws.Range("E10") '  +1.28%  '

# Row 11
Set-TextValue $ws.Range("D11") '5.79'
Set-TextValue $ws.Range("E11") '  -3.04%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.501'
Set-TextValue $ws.Range("E12") '  -2.37%  '

# Row 13
Set-TextValue $ws.Range("D13") '0.0000271'

# Row 14
Set-TextValue $ws.Range("D14") '38.84'
Set-TextValue $ws.Range("E14") '  -0.56%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.769.15'
Set-TextValue $ws.Range("E15") '  +1.25%  '

# Row 16
Set-TextValue $ws.Range("D16") '66.633.90'
Set-TextValue $ws.Range("E16") '  -0.12%  '

# Row 17
Set-TextValue $ws.Range("B17") 'WrappedEther'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D17") '3.239.06'
Set-TextValue $ws.Range("E17") '  +1.29%  '

# Row 18
Set-TextValue $ws.Range("B18") 'Polkadot'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D18") '7.28'
Set-TextValue $ws.Range("E18") '  -2.14%  '

# Row 19
Set-TextValue $ws.Range("B19") 'TRON'
Set-TextValue $ws.Range("C19") 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D19") '0.113'
Set-TextValue $ws.Range("E19") '  +1.30%  '

# Row 20
Set-TextValue $ws.Range("D20") '507.75'
Set-TextValue $ws.Range("E20") '  -1.71%  '

# Row 21
Set-TextValue $ws.Range("D21") '15.22'
Set-TextValue $ws.Range("E21") '  -0.99%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.740'
Set-TextValue $ws.Range("E22") '  +0.06%  '

# Row 23
Set-TextValue $ws.Range("E23") '  -2.26%  '

# Row 24
Set-TextValue $ws.Range("E24") '  -2.88%  '

# Row 25
Set-TextValue $ws.Range("D25") '86.19'
Set-TextValue $ws.Range("E25") '  +1.22%  '

# Row 26
Set-TextValue $ws.Range("D26") '0.167'
Set-TextValue $ws.Range("E26") '  +86.00%  '

# Row 27
Set-TextValue $ws.Range("E27") '  +0.04%  '

# Row 28
Set-TextValue $ws.Range("D28") '3.00'
Set-TextValue $ws.Range("E28") '  -0.67%  '

# Row 29
Set-TextValue $ws.Range("D29") '9.02'
Set-TextValue $ws.Range("E29") '  -3.02%  '

# Row 30
Set-TextValue $ws.Range("E30") '  -3.14%  '

# Row 31
Set-TextValue $ws.Range("B31") 'Stacks'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D31") '2.91'
Set-TextValue $ws.Range("E31") '  -6.29%  '

# Row 32
Set-TextValue $ws.Range("B32") 'NEARProtocol'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D32") '6.93'
Set-TextValue $ws.Range("E32") '  -1.90%  '

# Row 33
Set-TextValue $ws.Range("D33") '28.18'
Set-TextValue $ws.Range("E33") '  -0.08%  '

# Row 34
Set-TextValue $ws.Range("E34") '  +0.09%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -5.13%  '

# Row 36
Set-TextValue $ws.Range("E36") '  -3.53%  '

# Row 37
Set-TextValue $ws.Range("B37") 'PEPE'
Set-TextValue $ws.Range("C37") 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D37") '0.0₃0793'
Set-TextValue $ws.Range("E37") '  +14.47%  '

# Row 38
Set-TextValue $ws.Range("B38") 'OKB'
Set-TextValue $ws.Range("C38") 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D38") '55.34'
Set-TextValue $ws.Range("E38") '  +0.78%  '

# Row 39
Set-TextValue $ws.Range("D39") '494.02'
Set-TextValue $ws.Range("E39") '  -4.66%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +7.46%  '

# Row 41
Set-TextValue $ws.Range("E41") '  -0.96%  '

# Row 42
Set-TextValue $ws.Range("E42") '  +2.43%  '

# Row 43
Set-TextValue $ws.Range("D43") '8.70'
Set-TextValue $ws.Range("E43") '  -2.40%  '

# Row 44
Set-TextValue $ws.Range("E44") '  -4.57%  '

# Row 45
Set-TextValue $ws.Range("D45") '2.949.01'
Set-TextValue $ws.Range("E45") '  +2.46%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -1.55%  '

# Row 47
Set-TextValue $ws.Range("D47") '28.12'
Set-TextValue $ws.Range("E47") '  -2.16%  '

# Row 48
Set-TextValue $ws.Range("D48") '2.41'
Set-TextValue $ws.Range("E48") '  -0.43%  '

# Row 49
Set-TextValue $ws.Range("E49") '  +0.92%  '

# Row 50
Set-TextValue $ws.Range("E50") '  -0.05%  '

# Row 51
Set-TextValue $ws.Range("E51") '  -3.10%  '
